$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Insert 3 new rows at row 13 for the new "Omega SouthShore" SharePoint asset entries.
$ws.Rows.Item(13).Resize(3).Insert() | Out-Null

$ws.Cells.Item(13,1).Value2 = "TPL_Omega_SouthShore_SharePoint_url"
$ws.Cells.Item(13,2).Value2 = "TPL_Omega_SouthShore_SharePoint_url"

$ws.Cells.Item(14,1).Value2 = "TPL_Omega_SouthShore_SharePoint_LibraryName"
$ws.Cells.Item(14,2).Value2 = "TPL_Omega_SouthShore_SharePoint_LibraryName"

$ws.Cells.Item(15,1).Value2 = "TPL_Omega_SouthShore_SharePoint_FolderName"
$ws.Cells.Item(15,2).Value2 = "TPL_Omega_SouthShore_SharePoint_FolderName"

# The former "...ShareDrive_Path" rows (now shifted to rows 16-19) are renamed to
# "...ShareDrive_Folder" to reflect the new naming convention.
$ws.Cells.Item(16,1).Value2 = "TPL_eCWProductivity_Report_ShareDrive_Folder"
$ws.Cells.Item(16,2).Value2 = "TPL_eCW_Productivity_Report_ShareDrive_Folder"

$ws.Cells.Item(17,1).Value2 = "TPL_OmegaNYP_Report_ShareDrive_Folder"
$ws.Cells.Item(17,2).Value2 = "TPL_omega_NYP_Report_ShareDrive_Folder"

$ws.Cells.Item(18,1).Value2 = "TPL_OmegaUCC_Report_ShareDrive_Folder"
$ws.Cells.Item(18,2).Value2 = "TPL_omega_UCC_Report_ShareDrive_Folder"

$ws.Cells.Item(19,1).Value2 = "TPL_GebbsReport_ShareDrive_Folder"
$ws.Cells.Item(19,2).Value2 = "TPL_Gebbs_Report_ShareDrive_Folder"

# Insert 2 new rows right after (before the former row 17 "TPL_DataBaseName", now row 20)
# for the new South Shore share-drive folder asset and the generic share-drive URL asset.
$ws.Rows.Item(20).Resize(2).Insert() | Out-Null

$ws.Cells.Item(20,1).Value2 = "TPL_Omega_SouthShore_ShareDrive_Folder"
$ws.Cells.Item(20,2).Value2 = "TPL_Omega_SouthShore_ShareDrive_Folder"

$ws.Cells.Item(21,1).Value2 = "TPL_ShareDrive_URL"
$ws.Cells.Item(21,2).Value2 = "TPL_ShareDrive_URL"

# Update the sheet selection to match the new workbook state.
$ws.Activate()
$ws.Range("K12").Select() | Out-Null
